# edit.ps1 - apply "Se agregaron pasos al documento de detalle" changes
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert two new paragraphs right after the first paragraph
#    ("Pasos de los puntos") and before the "Para el firewall:" paragraph.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$pNew1 = $d.Paragraphs.Item(2)
$pNew1.Range.Text = "A través del VMware se indicaron 2 redes diferentes como internas para las del grupo 192.168.10.x y 192.168.20.x"

$pNew1.Range.InsertParagraphAfter()
$pNew2 = $d.Paragraphs.Item(3)
$pNew2.Range.Text = "Se configuraron los IP estáticos en las máquinas según se indica en el gráfico."

# ---------------------------------------------------------------------------
# Helper: wrap a <w:p>...</w:p> fragment in the WordProcessingML package
# envelope expected by Range.InsertXML, then replace the given paragraph's
# Range contents with it (InsertXML replaces the target range in place).
# ---------------------------------------------------------------------------
function Set-ParagraphXml($paragraph, [string]$fragment) {
    $envelope = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $fragment +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $paragraph.Range.InsertXML($envelope)
}

# ---------------------------------------------------------------------------
# 2) "Se creo el archivo rules.ipv4.conf ..." paragraph: wrap "rules.ipv4.conf"
#    in gramStart/gramEnd proofErr markers (now at paragraph index 4, after
#    the two paragraphs inserted above).
# ---------------------------------------------------------------------------
$firewallFragment = @'
<w:p><w:r><w:br/></w:r><w:r><w:t>Para el firewall:</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Se creo el archivo </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>rules.ip</w:t></w:r><w:r><w:t>v</w:t></w:r><w:r><w:t>4</w:t></w:r><w:r><w:t>.conf</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">para guardar las reglas de las </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>iptables</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:br/><w:t xml:space="preserve">Se indicó </w:t></w:r><w:r><w:t>en el archivo de carga de reglas</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>network</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>if-up.d</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>iptables</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">que haga </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>restore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de </w:t></w:r><w:r><w:t>este archivo creado</w:t></w:r></w:p>
'@

$firewallPara = $d.Paragraphs.Item(4)
Set-ParagraphXml $firewallPara $firewallFragment

# ---------------------------------------------------------------------------
# 3) "Se indicaron las variables de entorno JAVA_HOME..." paragraph: split the
#    "~/." run into "~" + "/." and wrap "/.bashrc" in gramStart/gramEnd.
# ---------------------------------------------------------------------------
$javaHomeFragment = @'
<w:p><w:r><w:t>Se creó el directorio /</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>opt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:br/><w:t xml:space="preserve">Se copiaron los archivos </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jdk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y apache en el directorio /</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>opt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y se descomprimieron</w:t></w:r><w:r><w:br/><w:t>Se indicaron las variables de entorno JAVA_HOME, JRE_HOME y CATALINA_HOME en el archivo ~</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>/.</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bashrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>
'@

$javaHomePara = $d.Paragraphs.Item(9)
Set-ParagraphXml $javaHomePara $javaHomeFragment

# ---------------------------------------------------------------------------
# 4) "Se agregó en el archivo ~/.bashrc ..." paragraph: same split/markers.
# ---------------------------------------------------------------------------
$agregoFragment = @'
<w:p><w:r><w:t>Se agregó en el archivo ~</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>/.</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bashrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> el agregado del JAVA_HOME al PATH</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Se copio el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sample.war</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> en el directorio </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>webapp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tomcat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:br/><w:t xml:space="preserve">Se asignó permiso de ejecución al archivo startup.sh del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tomcat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:br/><w:t>Se ejecutó e inició el servidor apache</w:t></w:r></w:p>
'@

$agregoPara = $d.Paragraphs.Item(10)
Set-ParagraphXml $agregoPara $agregoFragment
